# Update "paises.xlsx" data: refresh COVID case numbers and re-sort a
# handful of countries whose case counts changed their ranking order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header: refresh the "last updated" timestamp -------------------------
$ws.Range("A1").Value = "Datos actualizados a 28 de Agosto de 2020 a las 21:25"

# --- Helper to write a full data row (B:H) ---------------------------------
function Set-RowStats($row, $casosTotales, $nuevosCasos, $casosActivos, $recuperados, $casosCriticos, $muertesHoy, $muertes) {
    $ws.Cells.Item($row, 2).Value = $casosTotales
    $ws.Cells.Item($row, 3).Value = $nuevosCasos
    $ws.Cells.Item($row, 4).Value = $casosActivos
    $ws.Cells.Item($row, 5).Value = $recuperados
    $ws.Cells.Item($row, 6).Value = $casosCriticos
    $ws.Cells.Item($row, 7).Value = $muertesHoy
    $ws.Cells.Item($row, 8).Value = $muertes
}

# --- Countries whose stats updated but keep their row/rank -----------------
# Estados Unidos (row 4)
Set-RowStats 4 6076842 30208 3356504 2534897 0 645 185441

# India (row 6)
Set-RowStats 6 3461240 76665 2647538 750989 0 1019 62713

# Francia (row 20)
Set-RowStats 20 267077 7379 86177 150304 0 20 30596

# Alemania (row 23)
Set-RowStats 23 241888 1323 215495 17033 0 1 9360

# Canada (row 27)
Set-RowStats 27 127074 226 113039 4927 0 6 9108

# Haiti (row 100)
Set-RowStats 100 8161 10 5743 2217 0 1 201

# --- Countries that swapped rank with their neighbour -----------------------
# Cabo Verde moves above Eslovaquia (rows 119/120)
$ws.Cells.Item(119, 1).Value = "Cabo Verde"
Set-RowStats 119 3745 46 2807 900 0 0 38

$ws.Cells.Item(120, 1).Value = "Eslovaquia"
Set-RowStats 120 3728 102 2225 1470 0 0 33

# Siria moves above Sudan del Sur (rows 133/134)
$ws.Cells.Item(133, 1).Value = "Siria"
Set-RowStats 133 2563 59 584 1876 0 3 103

$ws.Cells.Item(134, 1).Value = "Sudan del Sur"
Set-RowStats 134 2518 4 1290 1181 0 0 47

# Bahamas moves above Sierra Leona and Yemen (rows 140/141/142)
$ws.Cells.Item(140, 1).Value = "Bahamas"
Set-RowStats 140 2020 97 742 1228 0 2 50

$ws.Cells.Item(141, 1).Value = "Sierra Leona"
Set-RowStats 141 2013 0 1581 362 0 0 70

$ws.Cells.Item(142, 1).Value = "Yemen"
Set-RowStats 142 1943 10 1107 273 0 1 563
